$wb = $excel.ActiveWorkbook

# --- Data change: every 1 in the Problem1 adjacency/weight matrix becomes 0.5 ---
$ws1 = $wb.Worksheets.Item("Problem1")
$ws1.Range("A2:J11").Replace(1, 0.5)

# --- View change: Problem1 becomes the active/selected sheet & cell ---
$ws1.Activate()
$ws1.Range("I17").Select()
